# GHC7-waypoints.xlsx update: map/waypoints refresh.
#
# Changes applied:
#   1. The "Meter-to-Inch ratio" header/value pair in G1:H2 was swapped:
#        - H1 used to hold the bare ratio number (0.2) with no header above
#          it; it now carries the "Meter-to-Inch ratio" label.
#        - G2 used to (mistakenly) hold the "Meter-to-Inch ratio" label
#          under the "Bubble Diameter (in)" column; it now holds the
#          numeric ratio value (0.2).
#   2. Every offset formula in columns D and E (rows 2-33) is rewritten to
#      reference $G$2 (the new home of the ratio value) instead of $H$1.
#   3. Waypoint 19's "connected-to" list (F20) gains a connection to 20,
#      becoming "15;20;23" instead of "15;23".
#   4. The active selection moves to H2 (and the sheet scrolls back so
#      column A is in view again).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# -- 1. Swap the ratio label/value between H1 and G2 --------------------
$ws.Range("H1").Value = "Meter-to-Inch ratio"
$ws.Range("G2").Value = 0.2

# -- 2. Re-point every D/E offset formula at $G$2 instead of $H$1 -------
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=B$r+(`$G`$2/2)"
    $ws.Cells.Item($r, 5).Formula = "=C$r+(`$G`$2/2)"
}

# -- 3. Waypoint 19 (row 20) now also connects to waypoint 20 -----------
$ws.Range("F20").Value = "15;20;23"

# -- 4. Leave the sheet scrolled to the top-left, selection on H2 -------
[void]$ws.Range("A1").Select()
[void]$ws.Range("H2").Select()
